$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITR input data")
$ws.Range("A24").Value = "test"
Write-Host "Done"
